$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Data" values for No HP Siswa, Alamat Rumah Siswa, Nomor Telepon, and Foto Diri
$ws.Range("F3").Value = 252
$ws.Range("F4").Value = 252
$ws.Range("B6").Value = 252
$ws.Range("F9").Value = 214

# Extend the Total formula for the Personally Identifiable column to include row 9 (Foto Diri)
$ws.Range("F10").Formula = "=SUM(F1:F9)"
